# Apply the edits described by the commit "sync for ep 3"
$wb = $excel.ActiveWorkbook

# 1. Rename the "Salries" sheet to "Salaries"
$salSheet = $wb.Worksheets.Item("Salries")
$salSheet.Name = "Salaries"

# 2. On the Salaries sheet, fix the L3/L6 helper formulas (D77^2 -> E77)
$sal = $wb.Worksheets.Item("Salaries")

# Update the formula text: y=3.63E-10*1.46E-11*x -> y=77024+3110x
$sal.Range("P3").Value = "y=77024+3110x"

$sal.Range("L3").Formula = "=(COUNT(Table1[Salary])*(E77))-((Table1[[#Totals],[Years Employed]]^2))"
$sal.Range("L6").Formula = "=(COUNT(Table1[Salary])*(E77))-(Table1[[#Totals],[Years Employed]]^2)"

# 4. Add the Intercept/Slope labels + formulas in column N
$sal.Range("N1").Value = "Intercept"
$sal.Range("N3").Formula = "=INTERCEPT(Table1[Salary],Table1[Years Employed])"
$sal.Range("N4").Value = "Slope"
$sal.Range("N6").Formula = "=SLOPE(Table1[Salary],Table1[Years Employed])"

# 5. Switch active sheet / selection to Salaries, set zoom, and select N6
$sal.Activate()
$sal.Application.ActiveWindow.Zoom = 150
$sal.Range("N6").Select()
